$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 67.34448032620639
$ws.Range("B3").Value = 0.8472844316881497
$ws.Range("B4").Value = 0.08909779639285949
$ws.Range("B5").Value = 0.2685218614894114
